# Update evidence for task B5-B7
$wb = $excel.ActiveWorkbook

# --- B5: fill in the real evidence hashes (was placeholder text) ---
$b5 = $wb.Worksheets.Item("B5")
$b5.Range("A2").Value = "38EF6D816E6023D81C3C18B7A5AB13C4AF90672FBB24A71FB84C866F839CD5B7"
$b5.Range("A3").Value = "559B4D3C0168B6604706DA75DD04A5A65A0A1F5151B6FD8B6E1CE7ADF835B6FC"

# --- B6: fill in the real evidence hashes and move the selection to A3 ---
$b6 = $wb.Worksheets.Item("B6")
$b6.Range("A2").Value = "E7AA928546213422290DB7A141A2EEB3D7CA80D9F291865CFA3B77D823FCEAFF"
$b6.Range("A3").Value = "47A6FDDD23F2E47095193ACD964DADA94FD0B6EC2E4B29A21A3D6E70DEF83BA5"
$b6.Activate()
[void]$b6.Range("A3").Select()

# --- B7: new sheet, cloned from B6's layout, with its own evidence hashes ---
$b6.Copy($null, $b6)
$b7 = $wb.Worksheets.Item($wb.Worksheets.Count)
$b7.Name = "B7"
$b7.Range("A2").Value = "DD7A8876EA276D3D3279D63F8FB0A53B15E5639B9F2471BC9802A0F7B6313EED"
$b7.Range("A3").Value = "E6C6C2E27F6EB9990A1855CF60376ED95B2AF6B1D08841BAF0FF92CA80815FD9"

# B7 becomes the active/selected tab, with A3 the selected cell
$b7.Activate()
[void]$b7.Range("A3").Select()

# --- B5's own selection also moves to A3 ---
$b5.Activate()
[void]$b5.Range("A3").Select()

# --- B2's tab is no longer the selected one (B7 is now active) ---
# (selecting/activating B7 above already clears B2's tabSelected flag)

# Leave B7 as the final active sheet/tab, matching activeTab="25"
$b7.Activate()
